$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cve_distrito -> distrito name
$distritoNombre = @{
    1 = "SAN LUIS RIO COLORADO"
    2 = "NOGALES"
    3 = "HERMOSILLO"
    4 = "GUAYMAS"
    5 = "HERMOSILLO"
    6 = "CD. OBREGON"
    7 = "CD. OBREGON"
}

# Step 1: fix up the existing 20 data rows (rows 2-21).
# Column J held a 4-digit "cve_distrito" (2601/2602/2603); it should just be
# the district number, and columns K (distrito), L (cve_municipio) and M
# (municipio) need to be populated to match.
for ($r = 2; $r -le 21; $r++) {
    $oldJ = $ws.Cells.Item($r, 10).Value2
    $d = [int]$oldJ - 2600
    $ws.Cells.Item($r, 10).Value = $d
    $ws.Cells.Item($r, 11).Value = $distritoNombre[$d]
    $ws.Cells.Item($r, 12).Value = $d
    $ws.Cells.Item($r, 13).Value = $distritoNombre[$d]
}

# Step 2: append 20 more affiliates (rows 22-41), duplicating the A-I data of
# rows 2-21 (shifting id by +20) and continuing the district sequence
# (4, 5, 6, 7) across the new rows, also carrying over estatus/notas (O/P).
$newDistrito = @(4,4,4,4, 5,5,5,5,5,5, 6,6,6,6,6, 7,7,7,7,7)

for ($i = 0; $i -lt 20; $i++) {
    $srcRow = 2 + $i
    $dstRow = 22 + $i
    $d = $newDistrito[$i]

    $ws.Cells.Item($dstRow, 1).Value = $ws.Cells.Item($srcRow, 1).Value2 + 20   # id
    $ws.Cells.Item($dstRow, 2).Value = $ws.Cells.Item($srcRow, 2).Value2        # nombre
    $ws.Cells.Item($dstRow, 3).Value = $ws.Cells.Item($srcRow, 3).Value2        # sexo
    $ws.Cells.Item($dstRow, 4).Value = $ws.Cells.Item($srcRow, 4).Value2        # edad
    $ws.Cells.Item($dstRow, 5).Value = $ws.Cells.Item($srcRow, 5).Value2        # clave_elector
    $ws.Cells.Item($dstRow, 6).Value = $ws.Cells.Item($srcRow, 6).Value2        # email
    $ws.Cells.Item($dstRow, 7).Value = $ws.Cells.Item($srcRow, 7).Value2        # telefono
    $ws.Cells.Item($dstRow, 8).Value = $ws.Cells.Item($srcRow, 8).Value2        # cve_estado
    $ws.Cells.Item($dstRow, 9).Value = $ws.Cells.Item($srcRow, 9).Value2        # estado
    $ws.Cells.Item($dstRow, 10).Value = $d                                      # cve_distrito
    $ws.Cells.Item($dstRow, 11).Value = $distritoNombre[$d]                     # distrito
    $ws.Cells.Item($dstRow, 12).Value = $d                                      # cve_municipio
    $ws.Cells.Item($dstRow, 13).Value = $distritoNombre[$d]                     # municipio
    $ws.Cells.Item($dstRow, 15).Value = $ws.Cells.Item($srcRow, 15).Value2      # estatus
    $ws.Cells.Item($dstRow, 16).Value = $ws.Cells.Item($srcRow, 16).Value2      # notas
}

$ws.Range("L75").Select()
